# Trade #73 closed at 2026-02-17 21:12:53 - unknown UNKNOWN +0.000%
#
# This applies:
#  - Summary sheet roll-up counters (Total Trades, Win Rate %)
#  - Strategy Status roll-up for the MarketMaking row (Trades, Win Rate %)
#  - Closes the existing open MarketMaking trade (early exit) in both the
#    "All Trades" ledger and the per-strategy "MarketMaking" ledger
#  - Appends a brand-new OPEN MarketMaking trade row to both ledgers

$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell while keeping it text even when the
# text looks like a date/time/number (Excel would otherwise silently coerce
# a literal like "2026-02-17" into a date serial). We briefly force a text
# number format, assign the value, then restore a plain "Normal" style so we
# don't leave a stray format behind.
function Set-TextValue($ws, $addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 101     # Total Trades
$wsSummary.Range("B9").Value = 47.52   # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (row 5 = MarketMaking)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D5").Value = 68    # Trades
$wsStatus.Range("G5").Value = 50    # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Close out the previously OPEN trade (Trade # 101, row 102): it exits early.
$wsAll.Range("G102").Value = 0.06            # Exit Price
$wsAll.Range("H102").Value = "CLOSED"        # Status
$wsAll.Range("K102").Value = 101.15          # Capital After
Set-TextValue $wsAll "L102" "early_exit"     # Exit Reason
$wsAll.Range("M102").Value = 0.14            # Duration (min)

# Append the new trade (Trade # 134) as a fresh OPEN row.
Set-TextValue $wsAll "B135" "2026-02-17"
Set-TextValue $wsAll "C135" "21:12:47"
$wsAll.Range("A135").Value = 134
Set-TextValue $wsAll "D135" "MarketMaking"
Set-TextValue $wsAll "E135" "DOWN"
$wsAll.Range("F135").Value = 0.06
$wsAll.Range("H135").Value = "OPEN"
$wsAll.Range("I135").Value = 0
$wsAll.Range("J135").Value = 0
$wsAll.Range("K135").Value = 101.1496151053151
$wsAll.Range("M135").Value = 0
$wsAll.Range("N135").Value = 0
$wsAll.Range("O135").Value = 0
$wsAll.Range("P135").Value = 0.6
Set-TextValue $wsAll "Q135" "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet (per-strategy ledger; same trades, different column
# layout: L=Entry Slippage, M=Exit Slippage, N=Confidence, O=Entry Reason,
# P=Exit Reason, Q=Duration)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Close out the previously OPEN trade (Trade # 101, row 69).
$wsMM.Range("G69").Value = 0.06              # Exit Price
$wsMM.Range("H69").Value = "CLOSED"          # Status
$wsMM.Range("K69").Value = 101.15            # Capital After
Set-TextValue $wsMM "P69" "early_exit"       # Exit Reason
$wsMM.Range("Q69").Value = 0.14              # Duration (min)

# Append the new trade (Trade # 134) as a fresh OPEN row.
Set-TextValue $wsMM "B102" "2026-02-17"
Set-TextValue $wsMM "C102" "21:12:47"
$wsMM.Range("A102").Value = 134
Set-TextValue $wsMM "D102" "MarketMaking"
Set-TextValue $wsMM "E102" "DOWN"
$wsMM.Range("F102").Value = 0.06
$wsMM.Range("H102").Value = "OPEN"
$wsMM.Range("I102").Value = 0
$wsMM.Range("J102").Value = 0
$wsMM.Range("K102").Value = 101.1496151053151
$wsMM.Range("L102").Value = 0
$wsMM.Range("M102").Value = 0
$wsMM.Range("N102").Value = 0.6
Set-TextValue $wsMM "O102" "Normal spread capture: 19600 bps"
$wsMM.Range("Q102").Value = 0
